# Update "想去人数" (interested-count) figures on the "展览" (sheet1)
# and "全部类型" (sheet4) worksheets, per the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 549
$wsExpo.Range("F4").Value = 180
$wsExpo.Range("F5").Value = 296
$wsExpo.Range("F6").Value = 398
$wsExpo.Range("F7").Value = 250
$wsExpo.Range("F8").Value = 2332
$wsExpo.Range("F9").Value = 388
$wsExpo.Range("F10").Value = 5848
$wsExpo.Range("F11").Value = 145
$wsExpo.Range("F13").Value = 5

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 549
$wsAll.Range("F5").Value = 180
$wsAll.Range("F6").Value = 296
$wsAll.Range("F7").Value = 398
$wsAll.Range("F8").Value = 250
$wsAll.Range("F11").Value = 2332
$wsAll.Range("F12").Value = 388
$wsAll.Range("F13").Value = 5848
$wsAll.Range("F14").Value = 145
$wsAll.Range("F17").Value = 5
